# LV_T2012_GiftLogGiftRequestProcessGiftSubmittedGiftRequestDetailVerifyEditAndDeleteFunctionality.xlsx
# "Merge - GiftLogs - Tcs. Pages & Test Data - 9th Oct 2025"
#
# The GiftLogs test-data workbook's "test user" names were swapped out:
#   Melissa Zatta -> Julie Carthane
#   Zain Sheikh   -> Aja Mount
# These appear on the Users sheet (StdUser / ComplianceUser columns) and the
# StdUser ("SubmittedFor") value used on the GiftLog sheet.

$wb = $excel.ActiveWorkbook

$wsUsers    = $wb.Worksheets.Item("Users")
$wsGiftLog  = $wb.Worksheets.Item("GiftLog")
$wsGiftEdit = $wb.Worksheets.Item("GiftEdit")

# --- Users sheet: swap the standard user and the compliance user ---
$wsUsers.Range("A2").Value = "Julie Carthane"
$wsUsers.Range("B2").Value = "Aja Mount"

# --- GiftLog sheet: the submitted-for user changes along with it ---
$wsGiftLog.Range("B2").Value = "Julie Carthane"
# Row now wraps to two lines at the new default row height.
$wsGiftLog.Rows.Item(2).RowHeight = 30
$wsGiftLog.Range("B2").Select()

# --- GiftEdit sheet: row height re-flowed too; cursor left parked elsewhere ---
$wsGiftEdit.Rows.Item(2).RowHeight = 60
$wsGiftEdit.Range("E15").Select()

# --- Users sheet becomes the active tab/selection last, as in the saved file ---
$wsUsers.Range("B2").Select()
